# Auto-generated edit script applying numeric updates to the Leve profit tables
# across all 8 job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 5276.6
$ws.Range("I76").Value = 5096
$ws.Range("K76").Value = 5096
$ws.Range("M76").Value = -4781
$ws.Range("H79").Value = 5276.6
$ws.Range("I79").Value = 5096
$ws.Range("K79").Value = 5096
$ws.Range("M79").Value = -4004
$ws.Range("H99").Value = 1667719.9
$ws.Range("I99").Value = 338.33334
$ws.Range("K99").Value = 1015.00002
$ws.Range("M99").Value = 482.9999799999999
$ws.Range("H125").Value = 2612
$ws.Range("I125").Value = 1900
$ws.Range("J125").Value = 4036
$ws.Range("K125").Value = 17100
$ws.Range("L125").Value = 36324
$ws.Range("M125").Value = -14640
$ws.Range("N125").Value = -41244
$ws.Range("H137").Value = 51545.6
$ws.Range("I137").Value = 1277
$ws.Range("J137").Value = 144901.58
$ws.Range("K137").Value = 3831
$ws.Range("L137").Value = 434704.74
$ws.Range("M137").Value = -1281
$ws.Range("N137").Value = -439804.74
$ws.Range("H138").Value = 3127.3972
$ws.Range("J138").Value = 2712.638
$ws.Range("L138").Value = 8137.914
$ws.Range("N138").Value = -18417.914
$ws.Range("H141").Value = 1219466.2
$ws.Range("I141").Value = 1648193
$ws.Range("K141").Value = 4944579
$ws.Range("M141").Value = -4939399

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3928.9
$ws.Range("I32").Value = 3942.8384
$ws.Range("J32").Value = 2549
$ws.Range("K32").Value = 3942.8384
$ws.Range("L32").Value = 2549
$ws.Range("M32").Value = -3655.8384
$ws.Range("N32").Value = -3123
$ws.Range("H45").Value = 8183281.5
$ws.Range("I45").Value = 18000940
$ws.Range("J45").Value = 1899
$ws.Range("K45").Value = 18000940
$ws.Range("L45").Value = 1899
$ws.Range("M45").Value = -18000563
$ws.Range("N45").Value = -2653
$ws.Range("H61").Value = 27916.906
$ws.Range("I61").Value = 46769.332
$ws.Range("J61").Value = 3678.0715
$ws.Range("K61").Value = 46769.332
$ws.Range("L61").Value = 3678.0715
$ws.Range("M61").Value = -46557.332
$ws.Range("N61").Value = -4102.0715
$ws.Range("H74").Value = 968.36584
$ws.Range("I74").Value = 789.2162
$ws.Range("J74").Value = 2625.5
$ws.Range("K74").Value = 789.2162
$ws.Range("L74").Value = 2625.5
$ws.Range("M74").Value = 84.78380000000004
$ws.Range("N74").Value = -4373.5
$ws.Range("H77").Value = 968.36584
$ws.Range("I77").Value = 789.2162
$ws.Range("J77").Value = 2625.5
$ws.Range("K77").Value = 3946.081
$ws.Range("L77").Value = 13127.5
$ws.Range("M77").Value = 421.9190000000003
$ws.Range("N77").Value = -21863.5
$ws.Range("H102").Value = 2398.6
$ws.Range("I102").Value = 2398.6
$ws.Range("K102").Value = 2398.6
$ws.Range("M102").Value = -776.5999999999999
$ws.Range("H122").Value = 1849.4445
$ws.Range("I122").Value = 1820.7142
$ws.Range("J122").Value = 1950
$ws.Range("K122").Value = 5462.142599999999
$ws.Range("L122").Value = 5850
$ws.Range("M122").Value = -3012.142599999999
$ws.Range("N122").Value = -10750
$ws.Range("H136").Value = 27916.906
$ws.Range("I136").Value = 46769.332
$ws.Range("J136").Value = 3678.0715
$ws.Range("K136").Value = 140307.996
$ws.Range("L136").Value = 11034.2145
$ws.Range("M136").Value = -137757.996
$ws.Range("N136").Value = -16134.2145

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 251168.62
$ws.Range("I86").Value = 1429.8
$ws.Range("J86").Value = 667400
$ws.Range("K86").Value = 1429.8
$ws.Range("L86").Value = 667400
$ws.Range("M86").Value = -306.8
$ws.Range("N86").Value = -669646
$ws.Range("H89").Value = 251168.62
$ws.Range("I89").Value = 1429.8
$ws.Range("J89").Value = 667400
$ws.Range("K89").Value = 7149
$ws.Range("L89").Value = 3337000
$ws.Range("M89").Value = -1533
$ws.Range("N89").Value = -3348232
$ws.Range("H134").Value = 4301.4224
$ws.Range("I134").Value = 4566.5884
$ws.Range("J134").Value = 3481.818
$ws.Range("K134").Value = 13699.7652
$ws.Range("L134").Value = 10445.454
$ws.Range("M134").Value = -11164.7652
$ws.Range("N134").Value = -15515.454

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H42").Value = 10120
$ws.Range("J42").Value = 10120
$ws.Range("L42").Value = 10120
$ws.Range("N42").Value = -11306
$ws.Range("H105").Value = 859.125
$ws.Range("I105").Value = 948
$ws.Range("J105").Value = 770.25
$ws.Range("K105").Value = 948
$ws.Range("L105").Value = 770.25
$ws.Range("M105").Value = 799
$ws.Range("N105").Value = -4264.25
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 493.5
$ws.Range("I92").Value = 494
$ws.Range("K92").Value = 1482
$ws.Range("M92").Value = -234
$ws.Range("H97").Value = 1148.8889
$ws.Range("J97").Value = 2933.3333
$ws.Range("L97").Value = 8799.999899999999
$ws.Range("N97").Value = -9791.999899999999
$ws.Range("H109").Value = 25003926
$ws.Range("I109").Value = 62500944
$ws.Range("K109").Value = 187502832
$ws.Range("M109").Value = -187501792
$ws.Range("H113").Value = 38962.035
$ws.Range("J113").Value = 1158
$ws.Range("L113").Value = 3474
$ws.Range("N113").Value = -7814
$ws.Range("H137").Value = 5038
$ws.Range("J137").Value = 5329.476
$ws.Range("L137").Value = 15988.428
$ws.Range("N137").Value = -26188.428

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1142.1428
$ws.Range("I113").Value = 799
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 799
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = 1371
$ws.Range("N113").Value = -6340
$ws.Range("H135").Value = 80030
$ws.Range("J135").Value = 80030
$ws.Range("L135").Value = 80030
$ws.Range("N135").Value = -90170

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 9379
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("H46").Value = 2056.8
$ws.Range("I46").Value = 1394.8334
$ws.Range("K46").Value = 1394.8334
$ws.Range("M46").Value = -1206.8334
$ws.Range("H136").Value = 1352.7222
$ws.Range("I136").Value = 910.3570999999999
$ws.Range("K136").Value = 2731.0713
$ws.Range("M136").Value = -181.0712999999996

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 7938.3335
$ws.Range("I96").Value = 865.6
$ws.Range("K96").Value = 865.6
$ws.Range("M96").Value = 507.4
$ws.Range("H126").Value = 3620.5557
$ws.Range("I126").Value = 3461.2
$ws.Range("K126").Value = 10383.6
$ws.Range("M126").Value = -7913.599999999999
$ws.Range("H136").Value = 20577886
$ws.Range("I136").Value = 37038270
$ws.Range("K136").Value = 111114810
$ws.Range("M136").Value = -111112260
